$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row for Player 11, who did not return a card (N/R)
$ws.Range("A12").Value = "Player 11"
$ws.Range("B12").Value = "N/R"
$ws.Range("C12").Value = "N/R"
$ws.Range("D12").Value = "N/R"

# Move the active selection to reflect where the user finished editing
$ws.Range("E12").Select()
